$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Build the new "Full Name *" column (G) out of the existing
#    "First Name *" (G) and "Last Name *" (H) columns.
$ws.Range("G1").Value = "Full Name *"

for ($r = 2; $r -le 7; $r++) {
    $firstName = $ws.Cells.Item($r, 7).Value2
    $lastName  = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 7).Value = "$firstName $lastName"
}

# Remove the now-redundant "Last Name *" column (H); everything to the
# right (old I:R) shifts one column left (new H:Q). This also moves the
# hyperlinked "Seller Signatory Emails" column from K to J.
$ws.Columns("H").Delete()

# 2. Rename "Offer Quantity *" header to "Quantity *"
$ws.Range("A1").Value = "Quantity *"

# 3. Re-create the hyperlinks so their anchors reflect the new layout
#    (Email column stays F, Seller Signatory Emails moves K -> J).
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), "mailto:emp3@myfirm.com")
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), "mailto:emp1@investor1.com")
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), "mailto:emp1@investor2.com")
$ws.Hyperlinks.Add($ws.Cells.Item(4, 10), "mailto:emp3@myfirm.com")
$ws.Hyperlinks.Add($ws.Cells.Item(6, 10), "mailto:emp1@investor1.com")
$ws.Hyperlinks.Add($ws.Cells.Item(7, 10), "mailto:emp1@investor2.com")

# Restore the plain "Hyperlink" cell style (Add() otherwise leaves behind
# a duplicated style record).
$ws.Cells.Item(4, 6).Style  = "Hyperlink"
$ws.Cells.Item(6, 6).Style  = "Hyperlink"
$ws.Cells.Item(7, 6).Style  = "Hyperlink"
$ws.Cells.Item(4, 10).Style = "Hyperlink"
$ws.Cells.Item(6, 10).Style = "Hyperlink"
$ws.Cells.Item(7, 10).Style = "Hyperlink"

# 4. Select the new "Full Name *" column, as in the edited workbook.
$ws.Range("G1:G7").Select()
